$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'97.985.14"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.51%  "

$ws.Range("D3").Value = "'3.394.42"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.99%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").Value = "'253.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.92%  "

$ws.Range("D6").Value = "'679.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.80%  "

$ws.Range("D7").Value = "'1.43"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -7.21%  "

$ws.Range("D8").Value = "'0.426"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -7.91%  "

$ws.Range("B9").Value = "USDC"
$ws.Range("C9").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D9").Value = "'1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.06%  "

$ws.Range("B10").Value = "Cardano"
$ws.Range("C10").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D10").Value = "'1.04"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.00%  "

$ws.Range("D11").Value = "'3.390.93"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.98%  "

$ws.Range("D12").Value = "'0.214"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.70%  "

$ws.Range("D13").Value = "'41.31"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.16%  "

$ws.Range("D14").Value = "'6.22"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +10.87%  "

$ws.Range("D15").Value = "'97.650.83"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.67%  "

$ws.Range("D16").Value = "'0.0000262"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.82%  "

$ws.Range("D17").Value = "'4.025.26"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.02%  "

$ws.Range("D18").Value = "'8.81"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +13.67%  "

$ws.Range("D19").Value = "'3.387.59"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.81%  "

$ws.Range("D20").Value = "'0.567"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +27.05%  "

$ws.Range("D21").Value = "'17.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.34%  "

$ws.Range("D22").Value = "'10.86"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.87%  "

$ws.Range("D23").Value = "'3.42"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.32%  "

$ws.Range("D24").Value = "'504.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.59%  "

$ws.Range("D25").Value = "'0.0000202"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.88%  "

$ws.Range("D26").Value = "'6.48"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.98%  "

$ws.Range("D27").Value = "'99.08"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.58%  "

$ws.Range("D28").Value = "'12.50"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.24%  "

$ws.Range("D29").Value = "'3.586.85"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.38%  "

$ws.Range("D30").Value = "'0.149"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.56%  "

$ws.Range("D31").Value = "'11.39"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.82%  "

$ws.Range("D32").Value = "'0.999"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.23%  "

$ws.Range("E33").Value = "  +2.77%  "

$ws.Range("D34").Value = "'2.60"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +22.41%  "

$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.56%  "

$ws.Range("D36").Value = "'0.565"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.48%  "

$ws.Range("D37").Value = "'29.26"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.24%  "

$ws.Range("D38").Value = "'1.50"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +11.40%  "

$ws.Range("D39").Value = "'7.84"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.40%  "

$ws.Range("D40").Value = "'525.13"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.77%  "

$ws.Range("D41").Value = "'0.152"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.92%  "

$ws.Range("E42").Value = "  +0.00%  "

$ws.Range("D43").Value = "'24.71"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.01%  "

$ws.Range("D44").Value = "'0.859"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.17%  "

$ws.Range("B45").Value = "MantraDAO"
$ws.Range("C45").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D45").Value = "'3.77"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.86%  "

$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "'0.0431"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.80%  "

$ws.Range("B47").Value = "Cosmos"
$ws.Range("C47").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D47").Value = "'8.85"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +10.09%  "

$ws.Range("D48").Value = "'1.72"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +12.34%  "

$ws.Range("D49").Value = "'5.70"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +10.50%  "

$ws.Range("D50").Value = "'55.90"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +11.10%  "

$ws.Range("D51").Value = "'3.17"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.73%  "
